$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.149.15'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '3.365.86'
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.14%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '3.363.08'
$ws.Range("E8").Value = '  +1.06%  '
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.13%  '
$ws.Range("E11").Value = '  +3.57%  '
$ws.Range("E12").Value = '  +3.21%  '
$ws.Range("D13").Value = '3.949.51'
$ws.Range("E13").Value = '  +1.63%  '
$ws.Range("E14").Value = '  +1.94%  '
$ws.Range("E15").Value = '  +3.61%  '
$ws.Range("D16").Value = '3.368.72'
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.53%  '
$ws.Range("D18").Value = '61.270.49'
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("E19").Value = '  +3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '379.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.567'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("D24").Value = '3.507.98'
$ws.Range("E24").Value = '  +1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("E27").Value = '  +10.48%  '
$ws.Range("E28").Value = '  +12.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.36%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.89%  '
$ws.Range("E32").Value = '  +3.02%  '
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("D35").Value = '3.403.20'
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '160.45'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0785'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.762'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.96'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.97%  '
$ws.Range("D51").Value = '2.328.66'
$ws.Range("E51").Value = '  +7.65%  '
